$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (values, then pick up the same formatting as the other
# header cells via a format-only paste from an existing header cell)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New column data (I = I0, J = IF)
$data = @(
    @(8, 8),
    @(9, 9),
    @(7, 7),
    @(9, 9),
    @(6, 6),
    @(6, 6),
    @(12, 12),
    @(6, 7),
    @(8, 9),
    @(8, 8),
    @(8, 9),
    @(9, 9),
    @(7, 9),
    @(6, 6),
    @(7, 7),
    @(7, 8),
    @(6, 6),
    @(5, 5),
    @(5, 6),
    @(7, 8),
    @(5, 5),
    @(7, 7),
    @(4, 4),
    @(5, 5)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
